$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# === Header text updates (Volume/Number, and report week dates) ===
$ws.Range("A8").Value = "Volume 32   Number  16"
$ws.Range("C9").Value = "Report Covering the Week  4/14/2025  Through  4/20/2025"


# --- Style/type changes: copy formatting+shared-string text from a stable donor cell, then set numeric value if needed ---
$ws.Range("C14").Copy($ws.Range("C15"))
$ws.Range("C14").Copy($ws.Range("D15"))
$ws.Range("E14").Copy($ws.Range("E15"))
$ws.Range("F14").Copy($ws.Range("D22"))
$ws.Range("D22").Value = 1
$ws.Range("L14").Copy($ws.Range("E22"))
$ws.Range("E22").Value = -100
$ws.Range("C14").Copy($ws.Range("C23"))
$ws.Range("C14").Copy($ws.Range("D23"))
$ws.Range("E14").Copy($ws.Range("E23"))
$ws.Range("C14").Copy($ws.Range("C27"))
$ws.Range("C14").Copy($ws.Range("D27"))
$ws.Range("E14").Copy($ws.Range("E27"))
$ws.Range("F14").Copy($ws.Range("D28"))
$ws.Range("D28").Value = 2
$ws.Range("L14").Copy($ws.Range("E28"))
$ws.Range("E28").Value = -100
$ws.Range("F14").Copy($ws.Range("G28"))
$ws.Range("G28").Value = 2
$ws.Range("L14").Copy($ws.Range("H28"))
$ws.Range("H28").Value = 150
$ws.Range("F14").Copy($ws.Range("D31"))
$ws.Range("D31").Value = 1
$ws.Range("L14").Copy($ws.Range("E31"))
$ws.Range("E31").Value = -100
$ws.Range("F14").Copy($ws.Range("G31"))
$ws.Range("G31").Value = 1
$ws.Range("L14").Copy($ws.Range("H31"))
$ws.Range("H31").Value = 0
$ws.Range("C14").Copy($ws.Range("D33"))
$ws.Range("E14").Copy($ws.Range("E33"))

# --- Plain value changes ---
$ws.Range("N15").Value = 0
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = 60
$ws.Range("I16").Value = 32
$ws.Range("J16").Value = 50
$ws.Range("K16").Value = -36
$ws.Range("L16").Value = -47.540983606557
$ws.Range("M16").Value = -28.888888888888
$ws.Range("N16").Value = -89.808917197452
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 700
$ws.Range("F17").Value = 25
$ws.Range("G17").Value = 11
$ws.Range("H17").Value = 127.272727272727
$ws.Range("I17").Value = 74
$ws.Range("J17").Value = 77
$ws.Range("K17").Value = -3.896103896103
$ws.Range("L17").Value = 8.823529411764
$ws.Range("M17").Value = 138.709677419355
$ws.Range("N17").Value = -18.681318681318
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 56
$ws.Range("J18").Value = 45
$ws.Range("K18").Value = 24.444444444444
$ws.Range("L18").Value = -9.677419354838
$ws.Range("M18").Value = -16.417910447761
$ws.Range("N18").Value = -92.134831460674
$ws.Range("C19").Value = 8
$ws.Range("E19").Value = -38.461538461538
$ws.Range("F19").Value = 38
$ws.Range("H19").Value = -19.148936170212
$ws.Range("I19").Value = 159
$ws.Range("J19").Value = 197
$ws.Range("K19").Value = -19.289340101522
$ws.Range("L19").Value = -26.046511627907
$ws.Range("M19").Value = 34.745762711864
$ws.Range("N19").Value = -63.532110091743
$ws.Range("C20").Value = 9
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = 80
$ws.Range("F20").Value = 22
$ws.Range("G20").Value = 23
$ws.Range("H20").Value = -4.347826086956
$ws.Range("I20").Value = 72
$ws.Range("J20").Value = 70
$ws.Range("K20").Value = 2.857142857142
$ws.Range("L20").Value = 24.137931034482
$ws.Range("M20").Value = 30.90909090909
$ws.Range("N20").Value = -93.519351935193
$ws.Range("C21").Value = 28
$ws.Range("D21").Value = 25
$ws.Range("E21").Value = 12
$ws.Range("F21").Value = 107
$ws.Range("G21").Value = 100
$ws.Range("H21").Value = 7
$ws.Range("I21").Value = 400
$ws.Range("J21").Value = 444
$ws.Range("K21").Value = -9.909909909909
$ws.Range("L21").Value = -14.893617021276
$ws.Range("M21").Value = 26.182965299684
$ws.Range("N21").Value = -85.01872659176
$ws.Range("J22").Value = 6
$ws.Range("F23").Value = 2
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = -33.333333333333
$ws.Range("L23").Value = 5.555555555555
$ws.Range("C24").Value = 23
$ws.Range("E24").Value = 0
$ws.Range("G24").Value = 74
$ws.Range("H24").Value = 22.972972972973
$ws.Range("I24").Value = 344
$ws.Range("J24").Value = 323
$ws.Range("K24").Value = 6.501547987616
$ws.Range("L24").Value = -8.510638297872
$ws.Range("M24").Value = 16.216216216216
$ws.Range("C25").Value = 8
$ws.Range("E25").Value = 14.285714285714
$ws.Range("F25").Value = 27
$ws.Range("G25").Value = 25
$ws.Range("H25").Value = 8
$ws.Range("I25").Value = 134
$ws.Range("J25").Value = 117
$ws.Range("K25").Value = 14.529914529914
$ws.Range("L25").Value = -25.139664804469
$ws.Range("C26").Value = 12
$ws.Range("D26").Value = 10
$ws.Range("E26").Value = 20
$ws.Range("F26").Value = 43
$ws.Range("G26").Value = 30
$ws.Range("H26").Value = 43.333333333333
$ws.Range("I26").Value = 149
$ws.Range("J26").Value = 124
$ws.Range("K26").Value = 20.16129032258
$ws.Range("L26").Value = 36.697247706422
$ws.Range("M26").Value = 8.759124087591
$ws.Range("F28").Value = 5
$ws.Range("J28").Value = 10
$ws.Range("K28").Value = 90
$ws.Range("L28").Value = 90
$ws.Range("J31").Value = 13
$ws.Range("K31").Value = -46.153846153846
